$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.155.52'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.838.08'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.44%  '
$ws.Range('D5').Value = '244.64'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('E6').Value = '  +1.64%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('D8').Value = '0.07507'
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = '0.2931'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').Value = '23.17'
$ws.Range('E10').Value = '  +3.68%  '
$ws.Range('D11').Value = '0.07745'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').Value = '1.836.53'
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('D13').Value = '4.996'
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').Value = '0.6699'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('D15').Value = '82.73'
$ws.Range('E15').Value = '  +0.47%  '
$ws.Range('D16').Value = '0.000009352'
$ws.Range('E16').Value = '  -5.08%  '
$ws.Range('D17').Value = '6.009'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '29.168.60'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '2.079.37'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').Value = '12.59'
$ws.Range('E20').Value = '  +2.33%  '
$ws.Range('D21').Value = '223.97'
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').Value = '160.53'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').Value = '0.1401'
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').Value = '8.516'
$ws.Range('E27').Value = '  +1.08%  '
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').Value = '1.507'
$ws.Range('E29').Value = '  +1.27%  '
$ws.Range('D30').Value = '0.05983'
$ws.Range('E30').Value = '  +15.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.160'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.02%  '
$ws.Range('D32').Value = '4.068'
$ws.Range('E32').Value = '  +1.25%  '
$ws.Range('D33').Value = '1.207'
$ws.Range('E33').Value = '  +0.99%  '
$ws.Range('D34').Value = '0.7488'
$ws.Range('E34').Value = '  +1.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.850'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.46%  '
$ws.Range('E36').Value = '  +0.48%  '
$ws.Range('D37').Value = '2.685'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.773'
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.229.02'
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').Value = '6.572'
$ws.Range('E41').Value = '  +3.92%  '
$ws.Range('D42').Value = '0.8936'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = '1.006'
$ws.Range('E43').Value = '  +0.49%  '
$ws.Range('D44').Value = '102.37'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('D45').Value = '0.08019'
$ws.Range('E45').Value = '  +18.33%  '
$ws.Range('D46').Value = '1.980.36'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').Value = '0.00000000125'
$ws.Range('E47').Value = '  +2.44%  '
$ws.Range('D48').Value = '65.85'
$ws.Range('E48').Value = '  +3.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5110'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').Value = '0.4075'
$ws.Range('E50').Value = '  +1.72%  '
$ws.Range('D51').Value = '9.018'
$ws.Range('E51').Value = '  +2.29%  '
